# Add a new BOM line for the antenna used by the DLPLIGHTCRAFTER-compatible
# Zooid receiver build (A24-HASM-450 2.4GHz antenna from Mouser / Digi-Key).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table1 currently spans A1:G27 - extend it by one row; Excel will grow the
# table range/autofilter and the sheet dimension automatically.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Populate the new row (row 28). Write the brand-new text values in the same
# order they are first introduced left-to-right as read in the source data
# (Supplier Part Number 2, Supplier Part Number 1, Description) so that the
# shared-string table indices line up, then fill in the remaining columns
# that reuse already-existing shared strings.
$ws.Range("G28").Value = "A24-HASM-450-ND"
$ws.Range("E28").Value = "888-A24-HASM-450"
$ws.Range("C28").Value = "Antenna 2.4GHZ"
$ws.Range("B28").Value = 1
$ws.Range("D28").Value = "Mouser"
$ws.Range("F28").Value = "Digi-Key"

# Leave the active selection on the quantity cell of the newly added row.
$ws.Range("B28").Select() | Out-Null

# Match the page setup (paper size / orientation) recorded for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
